$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2020" column (Q) values, mirroring the existing style of
# the adjacent "2019" column (P) for each row.
$ws.Range("Q4").Value = 2020
$ws.Range("Q5").Value = 99.3
$ws.Range("Q6").Value = 99.371420589467803
$ws.Range("Q7").Value = 99.319469393395053
$ws.Range("Q8").Value = 99.442213297634979
$ws.Range("Q9").Value = 98.766881972988841
$ws.Range("Q10").Value = 99.212798374809537
$ws.Range("Q11").Value = 99.799160124155549
$ws.Range("Q12").Value = 99.146991622239156
$ws.Range("Q13").Value = 99.538370126605429
$ws.Range("Q14").Value = 99.765563948945029

# Copy the formatting from column P into the new column Q so the new
# cells inherit the same number formats / styles as their neighbours.
$ws.Range("P4:P14").Copy()
$ws.Range("Q4:Q14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the active selection to match the post-edit view state.
$ws.Range("P7").Select()
